# Refresh the cryptos listing (Price / Volume(1h) columns, and the
# PEPE / PancakeSwap row swap) to match the latest GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must stay as text
# (matching the original inlineStr representation in the workbook).
$textCells = @("D5","D6","D11","D12","D14","D19","D23","D24","D25","D28","D31","D33","D36","D37","D40","D41","D42","D44","D45","D46","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "56.282.92"
$ws.Range("E2").Value = "  +3.74%  "
$ws.Range("D3").Value = "2.316.08"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "518.93"
$ws.Range("E5").Value = "  +4.68%  "
$ws.Range("D6").Value = "133.98"
$ws.Range("E6").Value = "  +4.13%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").Value = "2.339.29"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("E10").Value = "  +8.33%  "
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  +7.56%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "24.07"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("D15").Value = "2.730.18"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "56.404.16"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("E17").Value = "  +4.88%  "
$ws.Range("D18").Value = "2.339.73"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("D19").Value = "10.55"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("E21").Value = "  +6.12%  "
$ws.Range("E22").Value = "  +5.00%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "60.78"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "0.992"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  +7.11%  "
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").Value = "171.78"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +12.34%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("E30").Value = "  +6.66%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  +5.94%  "
$ws.Range("E32").Value = "  +5.20%  "
$ws.Range("D33").Value = "18.37"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").Value = "1.27"
$ws.Range("E36").Value = "  +5.99%  "
$ws.Range("D37").Value = "0.930"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +8.62%  "
$ws.Range("E39").Value = "  +9.16%  "
$ws.Range("D40").Value = "37.52"
$ws.Range("E40").Value = "  +4.63%  "
$ws.Range("D41").Value = "0.383"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "141.12"
$ws.Range("E42").Value = "  +13.20%  "
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("D44").Value = "280.90"
$ws.Range("E44").Value = "  +16.34%  "
$ws.Range("D45").Value = "5.11"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("D46").Value = "0.0511"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("D47").Value = "0.0929"
$ws.Range("E47").Value = "  +4.00%  "
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("E50").Value = "  +5.95%  "
$ws.Range("D51").Value = "16.97"
$ws.Range("E51").Value = "  +5.51%  "
